# Update the "Förändrad" (Changed) date column (C) for rows 2-17
# from serial date 45171 (2023-09-02) to 45172 (2023-09-03).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C17").Value = 45172
